$d = $word.ActiveDocument

# Remove the leading "Pengganti Visite/Poliklinik<tab>:<tab>" label runs that
# precede the "${pengganti_visite}" placeholder, leaving just the placeholder
# text in that paragraph (mirrors the XML diff which drops those w:r runs).
$found = $d.Content.Find.Execute("Pengganti Visite/Poliklinik^t:^t", $true, $true, $false, $false, $false, $true, 1, $false, "", 2)
